# DaySale_2025-05-29_00-00.xlsx update:
#  - row for "فرشه اسنان دكتور فريش كبار" (Dr Fresh toothbrush) had its
#    remaining-qty / total / transaction-count figures corrected
#  - a new low-stock row was inserted for
#    "فرشه اسنان شاين اب بلاس اطفال " (Shine-up Plus kids toothbrush)
#  - the generated-on timestamp in the footer was refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row at 46 (pushes old 46.."كالونا" and 47.."معجون..."
#    down to 47/48, and the totals/footer rows down to 49/50).
# ---------------------------------------------------------------------
$ws.Rows.Item(46).Insert()

# Copy the (now shifted) row 48 formatting - identical column layout/style
# as every other data row - onto the freshly inserted blank row 46 so the
# new row picks up the same styles/borders used throughout the table.
$ws.Range("A48:Q48").Copy()
$ws.Range("A46:Q46").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row 45 - "فرشه اسنان دكتور فريش كبار": corrected figures
# ---------------------------------------------------------------------
$ws.Range("H45").Value2 = "1:0"
$ws.Range("P45").Value2 = "15.0000"
$ws.Range("Q45").Value2 = "1:0"

# ---------------------------------------------------------------------
# 3) Row 46 (new) - "فرشه اسنان شاين اب بلاس اطفال "
# ---------------------------------------------------------------------
$ws.Range("A46").Value2 = 40
$ws.Range("C46").Value2 = "فرشه اسنان شاين اب بلاس اطفال "
$ws.Range("H46").Value2 = "10:0"
$ws.Range("L46").Value2 = "0"
$ws.Range("N46").Value2 = "15.00"
$ws.Range("P46").Value2 = "45.0000"
$ws.Range("Q46").Value2 = "3:0"
$ws.Rows.Item(46).RowHeight = 25.5

# ---------------------------------------------------------------------
# 4) Restore the row height on the totals row, which used to be row 48
#    (ht 24.75) and is now row 49 (ht 25.5 in the refreshed report).
# ---------------------------------------------------------------------
$ws.Rows.Item(49).RowHeight = 25.5

# ---------------------------------------------------------------------
# 5) Refresh the "generated on" timestamp shown in the footer (was row
#    49, now row 50 after the insert).
# ---------------------------------------------------------------------
$ws.Range("A50").Value2 = "Thursday, 29 May, 2025 6:05 PM"
